# Applies the commit "removed ER tags from non-ER templates and non-ER tags"
# to the Swate metadata template workbook.
#
# 1. Rename the metadata sheet from "SwateTemplateMetadata" to "isa_template".
# 2. Clear the ER block (ER / ER Term Accession Number / ER Term Source REF)
#    values on the metadata sheet, since this template does not target a
#    specific endpoint repository.

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("SwateTemplateMetadata")
$wsMeta.Name = "isa_template"

$wsMeta.Range("B8").ClearContents()
$wsMeta.Range("B9").ClearContents()
$wsMeta.Range("B10").ClearContents()
